$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-31 Wednesday", "2024-08-01 Thursday"),
    @("53÷6=", "75÷9="),
    @("82÷4=", "75÷6="),
    @("74÷7=", "28÷9="),
    @("79÷4=", "75÷7="),
    @("16÷4=", "66÷4="),
    @("33÷2=", "97÷5="),
    @("31÷6=", "49÷5="),
    @("26÷7=", "30÷9="),
    @("25÷3=", "84÷6="),
    @("40÷4=", "74÷4="),
    @("61÷2=", "35÷4="),
    @("79÷8=", "14÷9="),
    @("67÷5=", "43÷9="),
    @("72÷5=", "51÷3="),
    @("57÷8=", "69÷2="),
    @("92÷4=", "25÷9="),
    @("48÷6=", "67÷5="),
    @("17÷3=", "95÷8="),
    @("61÷5=", "77÷5="),
    @("83÷3=", "44÷8="),
    @("22÷4=", "99÷7="),
    @("28÷4=", "75÷2="),
    @("93÷8=", "24÷4="),
    @("30÷8=", "35÷2="),
    @("83÷4=", "32÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
